$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the formula in C3: was =B3, now =+B4*3 (B4 is text "N/A" -> #VALUE! error)
$ws.Range("C3").Formula = "=+B4*3"

# Move the active selection from C9 to C10
$ws.Range("C10").Select()
